# Append: 2025-10-27 06:37 JST
# Update the "取得日時" (retrieved-at) timestamp in column A for the
# existing rows on the "ランサーズ" sheet from 06:29:20 to 06:37:44.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-27 06:37:44"

for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
